# Generate Report for handback
# Refresh the "Latest Handoff Datetime" (col D) and "Latest Handback DateTime"
# (col G) for the a6da3a44... localized file entry (row 3) on both the
# zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D3").Value = "2016-01-07 14:16:41"
$zhcn.Range("G3").Value = "2016-01-07 14:17:25"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D3").Value = "2016-01-07 14:16:53"
$dede.Range("G3").Value = "2016-01-07 14:17:46"
